# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price (D) and Volume(1h) (E) columns are plain text in this sheet (not numbers),
# so Price cells are forced to text via NumberFormat "@" before assignment, then
# restored to the default "Normal" style so no stray number formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.621.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.585.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.573.42"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.624"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.60%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.202"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.82%  "

$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000311"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.155.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.652.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.596.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("E19").Value = "  -0.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "562.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +12.23%  "

$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("E23").Value = "  -5.48%  "

$ws.Range("E24").Value = "  +7.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("E27").Value = "  +3.49%  "

$ws.Range("E28").Value = "  +1.16%  "

$ws.Range("E29").Value = "  -1.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.56%  "

$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "566.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("E37").Value = "  +7.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.36%  "

$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0779"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.43%  "

$ws.Range("E41").Value = "  +0.69%  "

$ws.Range("E42").Value = "  -3.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.344.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.05%  "

$ws.Range("E46").Value = "  -0.40%  "

$ws.Range("E47").Value = "  +1.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.00%  "

$ws.Range("E49").Value = "  +1.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +23.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.09%  "
